$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 2.05
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 44656

$ws.Range("B6").Value = 2.25
$ws.Range("C6").Value = 44656

$ws.Range("C6").Select()
